# T29R09 / T30R09 1940 data fix
# The 1940Survey tab mistakenly contained two rows belonging to a different
# township (T30R09): Segment_id "30-09-36-S" and "30-09-35-S". Remove them
# and re-sort the remaining rows by Segment_id (column C), which is what the
# author did after cleaning the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1940Survey")

$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()

# Find (and delete, bottom-most first so row indices of earlier matches
# stay valid) every row whose Segment_id (column C) belongs to T30R09.
for ($r = $lastRow; $r -ge 2; $r--) {
    $segId = $ws.Cells.Item($r, 3).Value()
    if ($segId -eq "30-09-36-S" -or $segId -eq "30-09-35-S") {
        $ws.Rows($r).Delete()
    }
}

$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()

# Re-sort the remaining data (rows 2..lastRow) by Segment_id ascending.
$sortRange = $ws.Range("A2:L" + $lastRow)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C" + $lastRow))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Match the author's final selection/scroll position left after deleting the
# two now-empty trailing rows.
$ws.Range("A" + ($lastRow + 1) + ":XFD" + ($lastRow + 2)).Select()
